# Fruta / hortaliza, semanal
# Update weekly date/price records: dates and associated volume/price values
# were re-sorted among the rows (rows 3, 4 and 7 are unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44193
$ws.Range("M2").Value = 40
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 3000

# Row 5
$ws.Range("D5").Value = 44175
$ws.Range("M5").Value = 25
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 4000

# Row 6
$ws.Range("D6").Value = 44196
$ws.Range("M6").Value = 56

# Row 8
$ws.Range("D8").Value = 44907
$ws.Range("M8").Value = 45
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 25000
$ws.Range("S8").Value = 5000

# Row 9
$ws.Range("D9").Value = 44181
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 20000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 20000
$ws.Range("S9").Value = 4000

# Row 10
$ws.Range("D10").Value = 44188
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("S10").Value = 3000

# Row 11
$ws.Range("D11").Value = 44189

# Row 12
$ws.Range("D12").Value = 44186
$ws.Range("M12").Value = 40
$ws.Range("N12").Value = 15000
$ws.Range("O12").Value = 15000
$ws.Range("P12").Value = 15000
$ws.Range("S12").Value = 3000

# Row 13
$ws.Range("D13").Value = 44179
$ws.Range("M13").Value = 45
$ws.Range("N13").Value = 20000
$ws.Range("O13").Value = 20000
$ws.Range("P13").Value = 20000
$ws.Range("S13").Value = 4000
